$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "33.767.47"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +8.10%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.777.10"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +4.52%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "225.44"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.561"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.79%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "30.51"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.41%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "46.76"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.39%  "
$ws.Range("E10").Value = "  +4.01%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0667"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.80%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0923"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.33%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "2.032.00"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.55%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.777.92"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.44%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.626"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.37%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "33.740.28"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +8.18%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "10.05"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "4.18"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "68.57"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.30%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "251.93"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.00%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.0₃0739"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.12%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.28"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.42%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "4.18"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.19%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.14"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.13%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "159.45"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.42%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "16.50"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("E28").Value = "  +1.36%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.94"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.91%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.82"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.70%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0513"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.05%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.20"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("E34").Value = "  +4.99%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.86"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +7.29%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.485.95"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0185"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.51%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "83.16"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  +1.93%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("E43").Value = "  +4.07%  "
$ws.Range("E44").Value = "  +2.39%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0511"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("E46").Value = "  +3.23%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.928.64"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +5.05%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "5.73"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.82%  "
$ws.Range("E49").Value = "  -0.13%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "11.77"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +13.69%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "50.73"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.61%  "
